# Add "Charles-Eric Letendart" to the list of authors (after "... Bettuzzi Luc")
# and move the "_GoBack" bookmark from its old position (around "dégâts" in the
# tower description section) to the end of that newly-appended text.

$d = $word.ActiveDocument

# --- Step 1: locate the end of the author list (" Luc") and append the new
# author as a separate run, with a temporary trailing placeholder character
# ("X") that lets us anchor the relocated bookmark precisely before we trim
# it back off again.
$rng = $d.Content
$rng.Find.Execute(" Luc", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.Text = ", Charles-Eric LetendartX"

# Force the newly-typed text to live in its own run (rather than being
# silently coalesced back into the preceding " Luc" run) by nudging its
# font size away from, and then back to, the inherited value.
$rng.Font.Size = 11
$rng.Font.Size = 10

# --- Step 2: move the "_GoBack" bookmark here, right before the "X"
# placeholder (i.e. immediately after "Letendart").
$bmPos = $rng.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Step 3: drop the placeholder character; the bookmark (being a simple
# position marker) stays put, ending up exactly after "Letendart" and before
# the paragraph mark.
$delRange = $d.Range($bmPos, $bmPos + 1)
$delRange.Text = ""
